# Updates the crypto price/volume table on Sheet1 (columns D = Price,
# E = Volume(1h)) to the latest scraped values.
#
# Several "Price" values are plain decimal-looking strings (e.g. "6.78",
# "1.00", "20.60"). Assigning those directly to Range.Value would make
# Excel auto-convert them to numbers (losing the trailing zero / text
# formatting, e.g. "1.00" -> 1). To keep them as text - matching the
# original workbook, where every Price/Volume cell is stored as a string
# - we prefix the literal with a leading apostrophe (the standard Excel
# "force text" trick) and then reset the cell style back to "Normal" so
# the quote-prefix flag doesn't leave a stray style behind. Values that
# can never be parsed as a plain number (multiple '.' groups, or the
# subscript-3 glyph) are assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.560.45"
$ws.Range("E2").Value = "  -2.08%  "

$ws.Range("D3").Value = "2.581.54"
$ws.Range("E3").Value = "  -2.59%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'542.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.76%  "

$ws.Range("D6").Value = "'144.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.13%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "'0.581"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.26%  "

$ws.Range("D9").Value = "'6.78"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.57%  "

$ws.Range("E10").Value = "  -3.00%  "

$ws.Range("E11").Value = "  +3.65%  "

$ws.Range("D13").Value = "3.034.91"
$ws.Range("E13").Value = "  -2.97%  "

$ws.Range("D14").Value = "58.474.77"
$ws.Range("E14").Value = "  -2.11%  "

$ws.Range("D15").Value = "'20.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.92%  "

$ws.Range("D16").Value = "2.577.78"
$ws.Range("E16").Value = "  -5.51%  "

$ws.Range("E17").Value = "  -2.87%  "

$ws.Range("D18").Value = "'4.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.73%  "

$ws.Range("D19").Value = "'334.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.13%  "

$ws.Range("D20").Value = "'10.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.25%  "

$ws.Range("D21").Value = "'6.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.12%  "

$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D23").Value = "'66.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.44%  "

$ws.Range("D24").Value = "'0.424"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.73%  "

$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.54%  "

$ws.Range("E26").Value = "  -5.02%  "

$ws.Range("D27").Value = "'7.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.08%  "

$ws.Range("D28").Value = "0.0₃0742"
$ws.Range("E28").Value = "  -1.58%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("E30").Value = "  -1.29%  "

$ws.Range("D31").Value = "'5.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.91%  "

$ws.Range("D32").Value = "'152.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.70%  "

$ws.Range("E33").Value = "  -0.70%  "

$ws.Range("E34").Value = "  -2.89%  "

$ws.Range("D35").Value = "'0.848"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.76%  "

$ws.Range("E36").Value = "  -4.69%  "

$ws.Range("D37").Value = "'0.821"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.91%  "

$ws.Range("E38").Value = "  -3.22%  "

$ws.Range("D39").Value = "'3.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.62%  "

$ws.Range("D40").Value = "'278.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.12%  "

$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").Value = "'0.592"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.36%  "

$ws.Range("D43").Value = "'10.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.96%  "

$ws.Range("E44").Value = "  -0.89%  "

$ws.Range("D45").Value = "'0.0528"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.57%  "

$ws.Range("D46").Value = "'18.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.04%  "

$ws.Range("E47").Value = "  +0.06%  "

$ws.Range("D48").Value = "1.901.92"
$ws.Range("E48").Value = "  -4.08%  "

$ws.Range("D49").Value = "'17.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.82%  "

$ws.Range("E50").Value = "  -3.27%  "

$ws.Range("D51").Value = "'109.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.19%  "
